# -------------------------------------------------------------------------
# bnr-cvd-2023-table1.xlsx : refresh table with 2023 data
#   - add a title row ("Table 1. Annual Event Count by Year")
#   - simplify header wording (drop coded parentheticals)
#   - re-order the annual rows to descending year (2023 on top)
#   - drop the 2009 row (fifteen years -> fourteen years, 2010-2023)
#   - recompute the Total row so it excludes 2009
#   - refresh the "Prepared by" footer date
# -------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new blank row at the top for the table title.
#    (this also pushes every existing merged range down by one row, which
#    is exactly what the target layout needs)
$ws.Rows.Item(1).Insert()
$ws.Range("A1").Value = "Table 1. Annual Event Count by Year"

# 2) Clean up the header wording (row numbers below are POST-insert).
$ws.Range("B2").Value = "CVD Event Type"
$ws.Range("B4").Value = "Patient Sex"
$ws.Range("A6").Value = "CVD Event Year"

# 3) Capture the 15 year-data rows (2009..2023, currently rows 7-21) plus
#    the Total row (22) before reshuffling anything.
$data = $ws.Range("A7:J21").Value()

# row 1 of $data is 2009 -- drop it, keep rows 2..15 (2010..2023) and
# reverse them so 2023 ends up on top.
$years = New-Object 'object[,]' 14,10
for ($i = 1; $i -le 14; $i++) {
    $srcRow = 16 - $i   # 15,14,...,2  => 2023 down to 2010
    for ($c = 1; $c -le 10; $c++) {
        $years[$i - 1, $c - 1] = $data[$srcRow, $c]
    }
}

# 4) Remove the now-unneeded 2009 row entirely (row 7, post-insert
#    numbering) -- everything below shifts up by one.
$ws.Rows.Item(7).Delete()

# 5) Write the reordered 2023->2010 block back into rows 7-20.
$ws.Range("A7:J20").Value = $years

# 6) Recompute the Total row (now row 21) from the visible data rows.
for ($c = 2; $c -le 10; $c++) {
    $colLetter = [char](64 + $c)
    $sum = 0
    for ($r = 7; $r -le 20; $r++) {
        $cellVal = $ws.Cells.Item($r, $c).Value()
        $sum = $sum + $cellVal
    }
    $ws.Cells.Item(21, $c).Value = $sum
}

# 7) Refresh the footer date.
$ws.Range("A22").Value = "Prepared by Ian Hambleton on 2025-11-22, for the Barbados National Registry"

"edit complete"
